# Auto-generated script applying numeric updates to Anima_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 100014130
$ws.Cells.Item(86, 9).Value = 228600960
$ws.Cells.Item(86, 10).Value = 2155.3333
$ws.Cells.Item(86, 11).Value = 228600960
$ws.Cells.Item(86, 12).Value = 2155.3333
$ws.Cells.Item(86, 13).Value = -228599837
$ws.Cells.Item(86, 14).Value = -4401.3333
$ws.Cells.Item(89, 8).Value = 100014130
$ws.Cells.Item(89, 9).Value = 228600960
$ws.Cells.Item(89, 10).Value = 2155.3333
$ws.Cells.Item(89, 11).Value = 1143004800
$ws.Cells.Item(89, 12).Value = 10776.6665
$ws.Cells.Item(89, 13).Value = -1142999184
$ws.Cells.Item(89, 14).Value = -22008.6665
$ws.Cells.Item(116, 8).Value = 2699.5
$ws.Cells.Item(116, 9).Value = 2758.75
$ws.Cells.Item(116, 10).Value = 2581
$ws.Cells.Item(116, 11).Value = 2758.75
$ws.Cells.Item(116, 12).Value = 2581
$ws.Cells.Item(116, 13).Value = 683.25
$ws.Cells.Item(116, 14).Value = -9465
$ws.Cells.Item(137, 8).Value = 3997.85
$ws.Cells.Item(137, 9).Value = 3264.111
$ws.Cells.Item(137, 11).Value = 9792.332999999999
$ws.Cells.Item(137, 13).Value = -7242.332999999999
$ws.Cells.Item(138, 8).Value = 1724.3158
$ws.Cells.Item(138, 9).Value = 1469.3549
$ws.Cells.Item(138, 10).Value = 1899.9556
$ws.Cells.Item(138, 11).Value = 4408.0647
$ws.Cells.Item(138, 12).Value = 5699.8668
$ws.Cells.Item(138, 13).Value = 731.9353000000001
$ws.Cells.Item(138, 14).Value = -15979.8668

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1315.625
$ws.Cells.Item(2, 9).Value = 1316.9286
$ws.Cells.Item(2, 10).Value = 1306.5
$ws.Cells.Item(2, 11).Value = 1316.9286
$ws.Cells.Item(2, 12).Value = 1306.5
$ws.Cells.Item(2, 13).Value = -1203.9286
$ws.Cells.Item(2, 14).Value = -1532.5
$ws.Cells.Item(32, 8).Value = 393380.75
$ws.Cells.Item(32, 9).Value = 433660.6
$ws.Cells.Item(32, 11).Value = 433660.6
$ws.Cells.Item(32, 13).Value = -433373.6
$ws.Cells.Item(45, 8).Value = 4134.294
$ws.Cells.Item(45, 9).Value = 4126
$ws.Cells.Item(45, 10).Value = 4154.2
$ws.Cells.Item(45, 11).Value = 4126
$ws.Cells.Item(45, 12).Value = 4154.2
$ws.Cells.Item(45, 13).Value = -3749
$ws.Cells.Item(45, 14).Value = -4908.2
$ws.Cells.Item(74, 8).Value = 876.88464
$ws.Cells.Item(74, 9).Value = 785.0952
$ws.Cells.Item(74, 11).Value = 785.0952
$ws.Cells.Item(74, 13).Value = 88.90480000000002
$ws.Cells.Item(77, 8).Value = 876.88464
$ws.Cells.Item(77, 9).Value = 785.0952
$ws.Cells.Item(77, 11).Value = 3925.476
$ws.Cells.Item(77, 13).Value = 442.5240000000003
$ws.Cells.Item(116, 8).Value = 1315.625
$ws.Cells.Item(116, 9).Value = 1316.9286
$ws.Cells.Item(116, 10).Value = 1306.5
$ws.Cells.Item(116, 11).Value = 1316.9286
$ws.Cells.Item(116, 12).Value = 1306.5
$ws.Cells.Item(116, 13).Value = 977.0714
$ws.Cells.Item(116, 14).Value = -5894.5
$ws.Cells.Item(122, 8).Value = 78131.16
$ws.Cells.Item(122, 9).Value = 112003.89
$ws.Cells.Item(122, 10).Value = 1917.5
$ws.Cells.Item(122, 11).Value = 336011.67
$ws.Cells.Item(122, 12).Value = 5752.5
$ws.Cells.Item(122, 13).Value = -333561.67
$ws.Cells.Item(122, 14).Value = -10652.5
$ws.Cells.Item(132, 8).Value = 3621.3547
$ws.Cells.Item(132, 9).Value = 2003.2
$ws.Cells.Item(132, 11).Value = 6009.6
$ws.Cells.Item(132, 13).Value = -3479.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1315.625
$ws.Cells.Item(3, 9).Value = 1316.9286
$ws.Cells.Item(3, 10).Value = 1306.5
$ws.Cells.Item(3, 11).Value = 1316.9286
$ws.Cells.Item(3, 12).Value = 1306.5
$ws.Cells.Item(3, 13).Value = -1202.9286
$ws.Cells.Item(3, 14).Value = -1534.5
$ws.Cells.Item(107, 8).Value = 49211.047
$ws.Cells.Item(107, 9).Value = 60278.35
$ws.Cells.Item(107, 11).Value = 60278.35
$ws.Cells.Item(107, 13).Value = -58358.35

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(100, 8).Value = 45000
$ws.Cells.Item(100, 10).Value = 45000
$ws.Cells.Item(100, 12).Value = 45000
$ws.Cells.Item(100, 14).Value = -47164

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(76, 8).Value = 3141.2856
$ws.Cells.Item(76, 10).Value = 3000
$ws.Cells.Item(76, 12).Value = 9000
$ws.Cells.Item(76, 14).Value = -9766
$ws.Cells.Item(79, 8).Value = 3141.2856
$ws.Cells.Item(79, 10).Value = 3000
$ws.Cells.Item(79, 12).Value = 9000
$ws.Cells.Item(79, 14).Value = -11652
$ws.Cells.Item(86, 8).Value = 763.3333
$ws.Cells.Item(86, 9).Value = 300
$ws.Cells.Item(86, 10).Value = 1226.6666
$ws.Cells.Item(86, 11).Value = 900
$ws.Cells.Item(86, 12).Value = 3679.9998
$ws.Cells.Item(86, 13).Value = 286
$ws.Cells.Item(86, 14).Value = -6051.9998
$ws.Cells.Item(89, 8).Value = 763.3333
$ws.Cells.Item(89, 9).Value = 300
$ws.Cells.Item(89, 10).Value = 1226.6666
$ws.Cells.Item(89, 11).Value = 2700
$ws.Cells.Item(89, 12).Value = 11039.9994
$ws.Cells.Item(89, 13).Value = 3228
$ws.Cells.Item(89, 14).Value = -22895.9994
$ws.Cells.Item(132, 8).Value = 1709.75
$ws.Cells.Item(132, 10).Value = 1935.875
$ws.Cells.Item(132, 12).Value = 17422.875
$ws.Cells.Item(132, 14).Value = -22482.875
$ws.Cells.Item(133, 8).Value = 13200.429
$ws.Cells.Item(133, 9).Value = 6460
$ws.Cells.Item(133, 11).Value = 19380
$ws.Cells.Item(133, 13).Value = -14320
$ws.Cells.Item(134, 8).Value = 6237.7295
$ws.Cells.Item(134, 9).Value = 3099.7273
$ws.Cells.Item(134, 11).Value = 9299.1819
$ws.Cells.Item(134, 13).Value = -4229.1819
$ws.Cells.Item(139, 8).Value = 3890.0557
$ws.Cells.Item(139, 9).Value = 1165.5
$ws.Cells.Item(139, 11).Value = 3496.5
$ws.Cells.Item(139, 13).Value = 1643.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 39187840
$ws.Cells.Item(80, 9).Value = 63628740
$ws.Cells.Item(80, 10).Value = 82400
$ws.Cells.Item(80, 11).Value = 63628740
$ws.Cells.Item(80, 12).Value = 82400
$ws.Cells.Item(80, 13).Value = -63627742
$ws.Cells.Item(80, 14).Value = -84396
$ws.Cells.Item(83, 8).Value = 39187840
$ws.Cells.Item(83, 9).Value = 63628740
$ws.Cells.Item(83, 10).Value = 82400
$ws.Cells.Item(83, 11).Value = 318143700
$ws.Cells.Item(83, 12).Value = 412000
$ws.Cells.Item(83, 13).Value = -318138708
$ws.Cells.Item(83, 14).Value = -421984
$ws.Cells.Item(97, 8).Value = 1576.6666
$ws.Cells.Item(97, 9).Value = 1318.5714
$ws.Cells.Item(97, 11).Value = 1318.5714
$ws.Cells.Item(97, 13).Value = -822.5714
$ws.Cells.Item(122, 8).Value = 1618.619
$ws.Cells.Item(122, 9).Value = 1665.7778
$ws.Cells.Item(122, 10).Value = 1335.6666
$ws.Cells.Item(122, 11).Value = 4997.3334
$ws.Cells.Item(122, 12).Value = 4006.9998
$ws.Cells.Item(122, 13).Value = -2547.3334
$ws.Cells.Item(122, 14).Value = -8906.9998
$ws.Cells.Item(132, 8).Value = 2608.3513
$ws.Cells.Item(132, 9).Value = 2370.3684
$ws.Cells.Item(132, 10).Value = 2859.5557
$ws.Cells.Item(132, 11).Value = 7111.1052
$ws.Cells.Item(132, 12).Value = 8578.667099999999
$ws.Cells.Item(132, 13).Value = -4581.1052
$ws.Cells.Item(132, 14).Value = -13638.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 125003250
$ws.Cells.Item(82, 9).Value = 166669330
$ws.Cells.Item(82, 11).Value = 166669330
$ws.Cells.Item(82, 13).Value = -166668969
$ws.Cells.Item(85, 8).Value = 125003250
$ws.Cells.Item(85, 9).Value = 166669330
$ws.Cells.Item(85, 11).Value = 166669330
$ws.Cells.Item(85, 13).Value = -166668082
$ws.Cells.Item(136, 8).Value = 10418546
$ws.Cells.Item(136, 9).Value = 1577.091
$ws.Cells.Item(136, 10).Value = 33335878
$ws.Cells.Item(136, 11).Value = 4731.272999999999
$ws.Cells.Item(136, 12).Value = 100007634
$ws.Cells.Item(136, 13).Value = -2181.272999999999
$ws.Cells.Item(136, 14).Value = -100012734

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4069.2856
$ws.Cells.Item(81, 9).Value = 4781.4287
$ws.Cells.Item(81, 10).Value = 3357.1428
$ws.Cells.Item(81, 11).Value = 9562.857400000001
$ws.Cells.Item(81, 12).Value = 6714.2856
$ws.Cells.Item(81, 13).Value = -8501.857400000001
$ws.Cells.Item(81, 14).Value = -8836.285599999999
$ws.Cells.Item(84, 8).Value = 4069.2856
$ws.Cells.Item(84, 9).Value = 4781.4287
$ws.Cells.Item(84, 10).Value = 3357.1428
$ws.Cells.Item(84, 11).Value = 47814.287
$ws.Cells.Item(84, 12).Value = 33571.428
$ws.Cells.Item(84, 13).Value = -42510.287
$ws.Cells.Item(84, 14).Value = -44179.428
$ws.Cells.Item(136, 8).Value = 2664.6843
$ws.Cells.Item(136, 9).Value = 2356.4583
$ws.Cells.Item(136, 10).Value = 3193.0715
$ws.Cells.Item(136, 11).Value = 7069.374899999999
$ws.Cells.Item(136, 12).Value = 9579.2145
$ws.Cells.Item(136, 13).Value = -4519.374899999999
$ws.Cells.Item(136, 14).Value = -14679.2145
